# KNX Klemmen Polaritaet vertauscht
#
# Applies the three/four changes captured in the commit:
#   1. Notes master "last saved" date field text (08.02.2020 -> 31.10.2020)
#   2. Merge the two "Jung "/"2138" runs in the device label into one run "Jung 2138"
#   3/4. Swap the position of the "-" and "+" KNX terminal labels

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) Notes master date placeholder field: 08.02.2020 -> 31.10.2020
# ---------------------------------------------------------------------------
try {
    $nm = $p.NotesMaster
    $dt = $nm.HeadersFooters.DateAndTime
    $dt.UseFormat = $false
    $dt.Value = "31.10.2020"
} catch {
}

try {
    $nm = $p.NotesMaster
    for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
        $phsh = $nm.Shapes.Item($i)
        if ($phsh.TextFrame.HasText -and $phsh.TextFrame.TextRange.Text -eq "08.02.2020") {
            $phsh.TextFrame.TextRange.Text = "31.10.2020"
        }
    }
} catch {
}

# ---------------------------------------------------------------------------
# locate shapes by their stable shape Id (survives any z-order/index drift)
# ---------------------------------------------------------------------------
function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $cand = $shapes.Item($i)
        if ($cand.Id -eq $id) {
            return $cand
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 2) "Jung " + "2138" runs -> single run "Jung 2138"
# ---------------------------------------------------------------------------
$lbl = Get-ShapeById $s.Shapes 3369992
$tr = $lbl.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i)
    $paraText = $para.Text.Trim()
    if ($paraText -eq "Jung 2138") {
        $full = $para.Characters(1, $paraText.Length)
        $full.Text = "Jung 2138"
    }
}

# ---------------------------------------------------------------------------
# 3/4) Swap the "-" / "+" KNX terminal label positions
# ---------------------------------------------------------------------------
$minus = Get-ShapeById $s.Shapes 172
$plus = Get-ShapeById $s.Shapes 173

$minus.Left = [double]"164.177490234375"
$minus.Top = [double]"5.937007904052734"

$plus.Left = [double]"156.39976501464844"
$plus.Top = [double]"6.9753546714782715"
